$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 1074
$ws1.Range("F7").Value = 1405
$ws1.Range("F12").Value = 131
$ws1.Range("F13").Value = 121
$ws1.Range("F14").Value = 414
$ws1.Range("F15").Value = 1305
$ws1.Range("F17").Value = 85
$ws1.Range("F20").Value = 637
$ws1.Range("F21").Value = 29
$ws1.Range("F22").Value = 196
$ws1.Range("F24").Value = 5603
$ws1.Range("F29").Value = 14130
$ws1.Range("F33").Value = 82
$ws1.Range("F34").Value = 419
$ws1.Range("F35").Value = 573
$ws1.Range("F36").Value = 4163
$ws1.Range("F37").Value = 105

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 1074
$ws4.Range("F7").Value = 1405
$ws4.Range("F12").Value = 131
$ws4.Range("F13").Value = 121
$ws4.Range("F14").Value = 414
$ws4.Range("F15").Value = 1305
$ws4.Range("F17").Value = 85
$ws4.Range("F21").Value = 637
$ws4.Range("F23").Value = 29
$ws4.Range("F24").Value = 196
$ws4.Range("F27").Value = 5603
$ws4.Range("F32").Value = 14130
$ws4.Range("F36").Value = 82
$ws4.Range("F37").Value = 419
$ws4.Range("F38").Value = 573
$ws4.Range("F39").Value = 4163
$ws4.Range("F40").Value = 105
